$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New worker/period roster written starting at row 16 (columns B:G)
# B = Tipo Doc Trabajador, C = N Doc Trabajador, D = Nombre Trabajador,
# E = Periodo Mora, F = Valor Mora, G = Salario Basico
$rows = @(
    @{R=16; C="1065578242"; D="JANER ALFONSO GUERRA FADUL";  E="1906"; F=33125; G=828116},
    @{R=17; C="1050552257"; D="JORGE ARMANDO MENDOZA CHOGO"; E="1906"; F=33125; G=828116},
    @{R=18; C="1002295442"; D="LINA YALILE CALVO MARTINEZ";  E="1906"; F=33125; G=828116},
    @{R=19; C="1065578242"; D="JANER ALFONSO GUERRA FADUL";  E="1907"; F=33125; G=828116},
    @{R=20; C="1050552257"; D="JORGE ARMANDO MENDOZA CHOGO"; E="1907"; F=33125; G=828116},
    @{R=21; C="1002295442"; D="LINA YALILE CALVO MARTINEZ";  E="1907"; F=33125; G=828116},
    @{R=22; C="1002295442"; D="LINA YALILE CALVO MARTINEZ";  E="1908"; F=33125; G=828116},
    @{R=23; C="1065578242"; D="JANER ALFONSO GUERRA FADUL";  E="1909"; F=33125; G=828116},
    @{R=24; C="1050552257"; D="JORGE ARMANDO MENDOZA CHOGO"; E="1909"; F=33125; G=828116},
    @{R=25; C="1002295442"; D="LINA YALILE CALVO MARTINEZ";  E="1909"; F=33125; G=828116},
    @{R=26; C="1065578242"; D="JANER ALFONSO GUERRA FADUL";  E="1910"; F=33125; G=828116},
    @{R=27; C="1050552257"; D="JORGE ARMANDO MENDOZA CHOGO"; E="1910"; F=33125; G=828116},
    @{R=28; C="1002295442"; D="LINA YALILE CALVO MARTINEZ";  E="1910"; F=33125; G=828116},
    @{R=29; C="1065578242"; D="JANER ALFONSO GUERRA FADUL";  E="1911"; F=33125; G=828116},
    @{R=30; C="1050552257"; D="JORGE ARMANDO MENDOZA CHOGO"; E="1911"; F=33125; G=828116},
    @{R=31; C="1002295442"; D="LINA YALILE CALVO MARTINEZ";  E="1911"; F=33125; G=828116},
    @{R=32; C="1065578242"; D="JANER ALFONSO GUERRA FADUL";  E="1912"; F=33125; G=828116},
    @{R=33; C="1050552257"; D="JORGE ARMANDO MENDOZA CHOGO"; E="1912"; F=33125; G=828116},
    @{R=34; C="1002295442"; D="LINA YALILE CALVO MARTINEZ";  E="1912"; F=33125; G=828116},
    @{R=35; C="1065578242"; D="JANER ALFONSO GUERRA FADUL";  E="2001"; F=23187; G=828116},
    @{R=36; C="1050552257"; D="JORGE ARMANDO MENDOZA CHOGO"; E="2001"; F=23187; G=828116},
    @{R=37; C="1002295442"; D="LINA YALILE CALVO MARTINEZ";  E="2001"; F=23187; G=828116}
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
}
